$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
@(0,"flowbot","raw",0,0.7271590955219772,0.4524416441885594,0.9308608766588276,0.774032934793638,1,0.4928359176784918,0.7674228387456954,0.8307695411224246,0,0.7405987079903645,0.8696467785633657,0.8425683631632674,0.8337775462800274,0.7950994895189494,0.8821729978304266,0.7177662161811016,0.317486956528631,0.9408167404984388,0.9575964388682972,0.930270759533416),
@(1,"dit","raw",0,0.5857348830409497,0.03254396573963395,0.4306949131800209,0.8326606359652551,1,0.0030892508486971,0.8688142912326473,0.839764567538702,0.643527804145209,0.8483753308821654,0.8143824953541777,0.8291440113688637,0.7611493274504246,0.4964071006064256,0.7559436593229824,0.07715669253684659,0.5207065632503656,0.9306116680171967,0.8720755800332262,0.9075981064118088),
@(2,"pndit","raw",0,0.6593728174787202,0.4534371586899295,0.9254487013554624,0.7165164737174714,0.8,0.4739295787849449,0.8432087025070387,0.6738864317119709,0.02910580919978817,0.6482108673969216,0.7654267278803587,0.818051281503638,0.7948428908696541,0.560030647584066,0.9319255824801966,0.6496204253872718,0.3587076775229521,0.9312246179022708,0.7659524601710614,0.9172368094373394),
@(3,"hisditonly","raw",0,0.4018885160030032,0.3999951530291327,0.2588553441903752,0.956465370369235,0.8,0.4552275937912256,0.9423531717800028,0.6167866004919149,0.6377611821860422,0.7888863170440094,0.8775635475843448,0.7874312500872318,0.7353826711300609,0.5349590702337547,0.7192495607071194,0,0.3168898335349684,0.9290759899563976,0.6243283866723438,0.5409737059831639),
@(4,"flowbot","sgp",0,0.7138278559596645,0.9289711072750856,0.9514460108732092,0.9518634328757024,1,0.9370209920996716,0.9242445997784808,0.8589363510185869,0.3082338336756623,0.933658158614168,0.9514086154617384,0.8948717248922784,0.8518262214254619,0.9263976798275305,0.9149756229122948,0.9257741993657416,0.3157468635106495,0.9262509398520498,0.9637726579785376,0.9240362455077578),
@(5,"dit","sgp",0,0.4888263669975369,0.4285037293703598,0.5999187990409752,0.9409981636998286,0.9993729480627284,0.9345694197791012,0.9404691822006718,0.8490230578053026,0.3057283838462095,0.931907603190636,0.9045179583939842,0.8692535253098979,0.7526504346555315,0.8555069896328563,0.4702669908623247,0.2080094774270662,0.6296906794596634,0.9405993364337688,0.7484414967927199,0.9185602400683186),
@(6,"pndit","sgp",0,0.6629126197714845,0.7466600982420376,0.8899383246470043,0.7415691147898343,1,0.9487199560367378,0.9108832142439263,0.8889656155154814,0.8401693753933762,0.8325445608477199,0.9144821807549848,0.8824768579847121,0.9044884399165648,0.844398905571329,0.5970163049040024,0.6507120977643042,0.3344129876238517,0.938724984809892,0.9602763930918902,0.9332042963984388),
@(7,"hisdit","sgp",0,0.7308486472314362,0.5122527009758682,0.8247648859714521,0.9307432024874372,1,0.9424919659479258,0.9539296841276552,0.6,0.9642704773391272,0.9313987455565628,0.9096867513268464,0.921167023838839,0.8219552328814616,0.7751844694860945,0.931173282872618,0.7614469929052529,0.6316867420732576,0.9430531400964764,0.8377400911706901,0.9308713802145232),
@(8,"pnhisdit","sgp",0,0.7064115310764051,0.5139790013335145,0.9253996850275632,0.774235605073848,0.9995528758264598,0.9605903252381828,0.9339458920182081,0.8079896561887102,0.3101927743762022,0.937363334411392,0.9631899269363112,0.9329073880727108,0.939052657691387,0.946467530764085,0.934873387371325,0.8968881266656368,0.817796311628736,0.9360559140643132,0.8537877913754456,0.9328740782906444),
@(9,"hisditonly","sgp",0,0.961556856628588,0.0889405835347683,0.2613828534507007,0.9410201566134628,1,0.9396101939917234,0.9589084133678036,0.8,0.6350873887600679,0.7902257641030861,0.932266467500007,0.8928541108167063,0.9037870681214863,0.9082900596185676,0.9229921306942644,0.0226352728866066,0.3537579683589634,0.93770971462465,0.6502681884210223,0.5771985225235223),
@(10,$null,$null,$null,4,2,4,11,5,2,11,5,3,15,49,153,9,5,3,1,3,8,7,4)
)

$nr = $rows.Count
$nc = $rows[0].Count
$arr = New-Object 'object[,]' $nr,$nc
for ($i = 0; $i -lt $nr; $i++) {
    for ($j = 0; $j -lt $nc; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A2:X12").Value = $arr

Write-Output "done"
